$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.522.62'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '3.660.92'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '628.87'
$ws.Range("E5").Value = '  -6.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.76'
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.12'
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.440'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000229'
$ws.Range("E12").Value = '  -2.97%  '
$ws.Range("D13").Value = '4.276.73'
$ws.Range("E13").Value = '  -1.18%  '
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '3.667.45'
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '69.530.93'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.87'
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.26'
$ws.Range("E20").Value = '  +4.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '469.46'
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.647'
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.67'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").Value = '3.805.09'
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000125'
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.05'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.65'
$ws.Range("E28").Value = '  -5.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.59'
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("E30").Value = '  -4.65%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.164'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.57'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("D36").Value = '3.660.97'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.31'
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '178.58'
$ws.Range("E39").Value = '  +3.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  -5.54%  '
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0888'
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.64'
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.94'
$ws.Range("E47").Value = '  -2.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.83'
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000263'
$ws.Range("E50").Value = '  -6.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.21'
$ws.Range("E51").Value = '  -5.77%  '